# "app: working TP plots"
# Adds a new "dimname" column (E) to the Tabelle1 sheet, giving each
# hypothesis row (n_par == 3, rows 10-16) the corresponding bain/restriktor
# dimension name, formats that new column with a slightly smaller font and
# vertical-centered alignment, bumps up the row height for those rows, and
# leaves the selection on B16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: header + values -----------------------------------
$ws.Range("E1").Value2  = "dimname"

$ws.Range("E10").Value2 = "H1.V1>V2>V3>0"
$ws.Range("E11").Value2 = "H1.complement"
$ws.Range("E12").Value2 = "H2.V1>V2>V3"
$ws.Range("E13").Value2 = "H2.complement"
$ws.Range("E14").Value2 = "H3.V1>V2&0.6666667*V1>V3"
$ws.Range("E15").Value2 = "H3.complement"
$ws.Range("E16").Value2 = "Hu"

# --- Formatting for the new data cells (E10:E16) ----------------------
# Build up the target style on E10 first (vertical-centered, 10pt Arial
# Unicode MS), then copy that exact format onto the rest of the column so
# every cell shares the same cell style.
$ws.Range("E10").VerticalAlignment = -4108   # xlCenter
$ws.Range("E10").Font.Size = 10
$ws.Range("E10").Font.Name = "Arial Unicode MS"

$ws.Range("E10").Copy() | Out-Null
$ws.Range("E11:E16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row heights for the newly-formatted rows --------------------------
$ws.Range("A10:E16").RowHeight = 15

# --- Restore the selection to match the saved view ---------------------
$ws.Range("B16").Select() | Out-Null
